$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292 (shifts existing rows 292-348 down to 293-349)
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new weekly record
$ws.Cells.Item(292, 1).Value  = 7
$ws.Cells.Item(292, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(292, 3).Value  = "Ñuble"
$ws.Cells.Item(292, 4).Value  = 45173
$ws.Cells.Item(292, 5).Value  = 16
$ws.Cells.Item(292, 6).Value  = 100112032
$ws.Cells.Item(292, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(292, 8).Value  = "Sin especificar"
$ws.Cells.Item(292, 9).Value  = "Primera"
$ws.Cells.Item(292, 10).Value = 100
$ws.Cells.Item(292, 11).Value = 14000
$ws.Cells.Item(292, 12).Value = 14000
$ws.Cells.Item(292, 13).Value = 14000
$ws.Cells.Item(292, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(292, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(292, 16).Value = 280
$ws.Cells.Item(292, 17).Value = 50
$ws.Cells.Item(292, 18).Value = "Hortaliza"
